# Update Name of Algo
# Applies updated KNN imputation values to column A, D, E of the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.257
$ws.Range("A6").Value = -22.291
$ws.Range("A7").Value = -20.029
$ws.Range("D7").Value = -8.154
$ws.Range("D12").Value = -7.302
$ws.Range("E13").Value = 16.561
$ws.Range("E14").Value = 17.328
$ws.Range("D15").Value = -8.196000000000002
$ws.Range("A16").Value = -21.887
$ws.Range("E16").Value = 16.767
$ws.Range("E19").Value = 16.57
$ws.Range("A20").Value = -20.14
$ws.Range("D20").Value = -7.714
$ws.Range("D21").Value = -7.988
$ws.Range("D22").Value = -8.025
$ws.Range("E22").Value = 16.752
$ws.Range("D23").Value = -7.997
$ws.Range("A28").Value = -22.111
$ws.Range("A29").Value = -21.469
$ws.Range("D29").Value = -7.400999999999999
$ws.Range("A32").Value = -21.648
$ws.Range("D34").Value = -7.904000000000001
$ws.Range("E36").Value = 16.901
$ws.Range("A40").Value = -19.896
$ws.Range("D42").Value = -7.958999999999999
$ws.Range("D43").Value = -7.853000000000002
$ws.Range("D44").Value = -7.831999999999999
$ws.Range("D45").Value = -7.531000000000001
$ws.Range("A46").Value = -21.707
$ws.Range("D46").Value = -7.784000000000001
$ws.Range("E46").Value = 17.031
$ws.Range("D50").Value = -8.250000000000002
$ws.Range("E50").Value = 16.473
$ws.Range("A51").Value = -22.069
$ws.Range("D51").Value = -8.280000000000001
$ws.Range("A52").Value = -22.16
$ws.Range("A57").Value = -22.263
$ws.Range("A59").Value = -22.312
$ws.Range("A62").Value = -22.085
$ws.Range("A66").Value = -21.532
$ws.Range("D66").Value = -7.557
$ws.Range("D67").Value = -7.145000000000001
$ws.Range("A73").Value = -20.213
$ws.Range("A74").Value = -21.244
$ws.Range("D79").Value = -7.642
$ws.Range("D84").Value = -8.300000000000001
$ws.Range("A92").Value = -21.582
$ws.Range("D92").Value = -6.597
$ws.Range("E95").Value = 17.238
$ws.Range("D97").Value = -8.259
$ws.Range("E97").Value = 16.766
$ws.Range("A100").Value = -22.217
